$d = $word.ActiveDocument

# 1. Rename the single {ExperiencePoint} placeholder to {ExperiencePoints[0]}
$d.Content.Find.Execute("{ExperiencePoint}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{ExperiencePoints[0]}", 2)

# 2. Append four more bullet paragraphs, cloning the formatting of the
#    (now last) ListParagraph-styled bullet, each holding the next
#    {ExperiencePoints[n]} placeholder.
for ($i = 1; $i -le 4; $i++) {
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $rng = $lastPara.Range
    $rng.Collapse(0)  # wdCollapseEnd
    $rng.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newPara.Range.Text = "{ExperiencePoints[$i]}"
}
